# Ajustando as porcentagens do relatório.
# Remove the two trailing empty cells on row 155 (E155/F155) and append four
# new attendance rows (156-159), mirroring the shape row 155 used to have
# (empty E/F placeholders) on the new last row (159).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 155 no longer carries the trailing empty E/F cells.
$ws.Range("E155:F155").ClearContents()

# Helper-ish literals: leading apostrophe forces text storage so values like
# "1", "2" and the yyyy-mm-dd looking dates are kept as literal text instead
# of being auto-coerced into numbers / date serials.

# Row 156
$ws.Range("A156").Value = "'1"
$ws.Range("B156").Value = "Pai do leonardo"
$ws.Range("C156").Value = "'2024-02-01"
$ws.Range("D156").Value = "P"

# Row 157
$ws.Range("A157").Value = "'1"
$ws.Range("B157").Value = "Pai do leonardo"
$ws.Range("C157").Value = "'2024-02-02"
$ws.Range("D157").Value = "A"

# Row 158
$ws.Range("A158").Value = "'2"
$ws.Range("B158").Value = "maria"
$ws.Range("C158").Value = "'2024-01-01"
$ws.Range("D158").Value = "P"

# Row 159 (also gets the trailing empty inline-string-style E/F cells that
# used to live on row 155)
$ws.Range("A159").Value = "'2"
$ws.Range("B159").Value = "maria"
$ws.Range("C159").Value = "'2024-01-02"
$ws.Range("D159").Value = "P"
$ws.Range("E159").Value = "'"
$ws.Range("F159").Value = "'"
